$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the header row (row 1) to the new, more descriptive column titles
# that document the expected data type / allowed values for each field.
$ws.Range("A1").Value = "Post Number (Integer)"
$ws.Range("B1").Value = "Publisher (String)"
$ws.Range("C1").Value = "Subject (Government, Violence, Health, Radicalism)"
$ws.Range("D1").Value = "Day (Integer, 1-8)"
$ws.Range("E1").Value = "Reaction (Happy, Sad, Angry)"
$ws.Range("F1").Value = "Hashtag 1 (String)"
$ws.Range("G1").Value = "Hashtag 2 (String)"
$ws.Range("H1").Value = "Hashtag 3 (String)"
$ws.Range("I1").Value = "Base Engagement (Integer)"
$ws.Range("J1").Value = "Boosted Engagement (Integer)"
$ws.Range("K1").Value = "Boost Cost (Float)"
$ws.Range("L1").Value = "Headline (String)"
$ws.Range("M1").Value = "Image File Path (String)"

# Resize the columns so the new, longer header text is fully visible.
$ws.Columns("A:M").AutoFit() | Out-Null

# Add a new syntax-warning note a few rows below the table.
$ws.Range("A6").Value = "Warning: do not use commas in string text"

# Leave the same cell selected as in the authored workbook.
$ws.Range("B12").Select() | Out-Null
